# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计", populated with
#    its fund-holding table (copied structure/style from the "2022-Q1"
#    sheet so fonts/borders/column layout match the other quarter sheets).
# 2. Prepend a corresponding row to the "总计" summary table and keep the
#    existing quarters below it (their running index in column A is
#    simply 0..4 top to bottom).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q3" sheet by copying "2022-Q1" (keeps header
# text + cell styles identical to its sibling quarter sheets), placed
# immediately after "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2022-Q1")
$templateSheet.Copy($null, $totalSheet)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# A "plain" (unstyled) helper cell used below to strip the quote-prefix
# style Excel applies when a numeric-looking string is forced to text.
$q3.Range("Z1").Value = "x"

# --- row 2 : 005571 ----------------------------------------------------
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'005571"
$q3.Range("C2").Value = "'中银证券新能源灵活配置混合A"
$q3.Range("D2").Value = "'0.53"
$q3.Range("E2").Value = "'90.32"
$q3.Range("F2").Value = "'4.79"
$q3.Range("G2").Value = "'0.0254"
$q3.Range("H2").Value = 10

# --- row 3 : 005572 ----------------------------------------------------
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'005572"
$q3.Range("C3").Value = "'中银证券新能源灵活配置混合C"
$q3.Range("D3").Value = "'0.25"
$q3.Range("E3").Value = "'90.32"
$q3.Range("F3").Value = "'4.79"
$q3.Range("G3").Value = "'0.0120"
$q3.Range("H3").Value = 10

# A3 needs the same bold/border style as A2 (copied sheet only carried
# the style for a single data row).
$q3.Range("A2").Copy()
$q3.Range("A3").PasteSpecial(-4122)

# Strip the quote-prefix style picked up by the textual numeric values
# above so the cells end up with the same plain style as the rest of
# the sheet (matches the other quarter sheets, whose text cells carry
# no explicit style).
$q3.Range("Z1").Copy()
$q3.Range("B2:G3").PasteSpecial(-4122)
$q3.Range("Z1").Clear()

$q3.Range("A1").Select()

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q3 row at the top of the "总计" data table,
# shifting the older quarters down by one row.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.04

# Restore the header-row style on the freshly inserted row's cells.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# Re-write the remaining rows explicitly (values only - keeps their
# existing style) so column A stays a clean 0..4 running index and the
# quarter figures line up with their correct row.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q1"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.08

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q4"
$totalSheet.Range("C4").Value = 18
$totalSheet.Range("D4").Value = 6.64

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q3"
$totalSheet.Range("C5").Value = 5
$totalSheet.Range("D5").Value = 0.57

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2020-Q4"
$totalSheet.Range("C6").Value = 2
$totalSheet.Range("D6").Value = 0.01

# The row inserted above picked up a bold/centered style from its
# neighbour - the data rows in this table carry no explicit style, so
# copy the (unstyled) format of row 3 onto the new row 2.
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)

# Leave the workbook on the summary sheet / first cell, matching the
# original file's active selection.
$totalSheet.Select()
$totalSheet.Range("A1").Select()
